$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.340.10'
$ws.Range("E2").Value = '  -0.08%  '
$ws.Range("D3").Value = '1.802.87'
$ws.Range("E3").Value = '  +0.86%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = '227.74'
$ws.Range("E5").Value = '  +0.68%  '
$ws.Range("D6").Value = '0.577'
$ws.Range("E6").Value = '  +3.94%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").Value = '36.16'
$ws.Range("E8").Value = '  +9.59%  '
$ws.Range("E9").Value = '  +2.18%  '
$ws.Range("E10").Value = '  +0.57%  '
$ws.Range("E11").Value = '  +2.09%  '
$ws.Range("D12").Value = '2.062.01'
$ws.Range("E12").Value = '  +0.78%  '
$ws.Range("D13").Value = '11.60'
$ws.Range("E13").Value = '  +3.96%  '
$ws.Range("D14").Value = '1.813.81'
$ws.Range("E14").Value = '  +1.43%  '
$ws.Range("E15").Value = '  +1.68%  '
$ws.Range("E16").Value = '  +4.88%  '
$ws.Range("D17").Value = '34.331.75'
$ws.Range("E17").Value = '  +0.04%  '
$ws.Range("D18").Value = '69.10'
$ws.Range("E18").Value = '  +0.99%  '
$ws.Range("D19").Value = '245.73'
$ws.Range("E19").Value = '  +0.20%  '
$ws.Range("E20").Value = '  +0.10%  '
$ws.Range("D21").Value = '11.51'
$ws.Range("E21").Value = '  +2.52%  '
$ws.Range("E22").Value = '  +0.08%  '
$ws.Range("D23").Value = '4.19'
$ws.Range("E23").Value = '  +0.90%  '
$ws.Range("D24").Value = '172.29'
$ws.Range("E24").Value = '  +2.22%  '
$ws.Range("D25").Value = '2.14'
$ws.Range("E25").Value = '  +3.47%  '
$ws.Range("D26").Value = '7.93'
$ws.Range("E26").Value = '  +8.30%  '
$ws.Range("E27").Value = '  +1.75%  '
$ws.Range("E28").Value = '  +2.92%  '
$ws.Range("D30").Value = '4.06'
$ws.Range("E30").Value = '  +1.07%  '
$ws.Range("E31").Value = '  +1.04%  '
$ws.Range("E32").Value = '  +1.31%  '
$ws.Range("D33").Value = '1.24'
$ws.Range("E33").Value = '  +1.09%  '
$ws.Range("D34").Value = '1.83'
$ws.Range("E34").Value = '  +0.46%  '
$ws.Range("D35").Value = '1.394.65'
$ws.Range("E35").Value = '  -1.02%  '
$ws.Range("D36").Value = '0.675'
$ws.Range("E36").Value = '  -1.05%  '
$ws.Range("D37").Value = '2.48'
$ws.Range("E37").Value = '  -5.23%  '
$ws.Range("E38").Value = '  -0.39%  '
$ws.Range("E39").Value = '  -0.11%  '
$ws.Range("D40").Value = '1.24'
$ws.Range("E40").Value = '  +11.54%  '
$ws.Range("D41").Value = '0.962'
$ws.Range("E41").Value = '  +2.51%  '
$ws.Range("E42").Value = '  +1.01%  '
$ws.Range("D43").Value = '82.02'
$ws.Range("E43").Value = '  -2.71%  '
$ws.Range("E44").Value = '  +0.40%  '
$ws.Range("D45").Value = '13.56'
$ws.Range("E45").Value = '  -3.83%  '
$ws.Range("D46").Value = '6.04'
$ws.Range("E47").Value = '  -5.15%  '
$ws.Range("D48").Value = '1.963.20'
$ws.Range("E48").Value = '  +0.88%  '
$ws.Range("D49").Value = '104.93'
$ws.Range("E49").Value = '  -0.50%  '
$ws.Range("E50").Value = '  +0.07%  '
$ws.Range("E51").Value = '  -0.30%  '
